# Generate Report for Handback
# - Marks the two localization files as handed back (Overview status text)
# - Fills in the "Latest Target File" / "Latest Handback File" / "Latest
#   Handback DateTime" columns on the zh-cn and de-de report sheets, adding
#   hyperlinks on the newly populated "Latest Target File" cells
# - Widens a few columns that now hold longer text

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Overview sheet: flip the per-language status from "Ready for
#    handoff" to "Handed back: in sync with en-US" and widen the status
#    columns (E = zh-cn, F = de-de) so the longer text fits.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Range("E:F").ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# Helper data: display name -> hyperlink target URL, read off the
# existing "Source File Name" hyperlinks (column A) so the newly added
# "Latest Target File" hyperlinks point at the same place.
# ---------------------------------------------------------------------
function Get-HyperlinkMap($ws) {
    $map = @{}
    foreach ($h in $ws.Hyperlinks) {
        $map[$h.TextToDisplay] = $h.Address
    }
    return $map
}

$fileA = "53393bc4-cd0b-4267-a893-b0bcde67eabe.md"
$fileB = "e075fa3d-c089-4848-bda2-6b19a07b1411.md"

# ---------------------------------------------------------------------
# 2. zh-cn report sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhLinks = Get-HyperlinkMap $wsZh

# Row 2 (53393bc4-...)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhLinks[$fileA], "", "", $fileA)
$wsZh.Range("J2").Value = "53393bc4-cd0b-4267-a893-b0bcde67eabe.3693b69075ba1d17482e18fce47f9862ff465be6.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-22 10:09:19"

# Row 3 (e075fa3d-...)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhLinks[$fileB], "", "", $fileB)
$wsZh.Range("J3").Value = "e075fa3d-c089-4848-bda2-6b19a07b1411.9dd91f8f50e743ff1c20550dc3d6ce6f44f5f606.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-22 10:09:19"

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Range("I:J").ColumnWidth = 40

# ---------------------------------------------------------------------
# 3. de-de report sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$deLinks = Get-HyperlinkMap $wsDe

# Row 2 (53393bc4-...)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deLinks[$fileA], "", "", $fileA)
$wsDe.Range("J2").Value = "53393bc4-cd0b-4267-a893-b0bcde67eabe.3693b69075ba1d17482e18fce47f9862ff465be6.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-22 10:09:26"

# Row 3 (e075fa3d-...)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deLinks[$fileB], "", "", $fileB)
$wsDe.Range("J3").Value = "e075fa3d-c089-4848-bda2-6b19a07b1411.9dd91f8f50e743ff1c20550dc3d6ce6f44f5f606.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-22 10:09:26"

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Range("I:J").ColumnWidth = 40
